# Update course excel file: replace the single "FACULTY OF BUSINESS &
# TECHNOLOGY" department label with per-row department names, clear the
# old per-row promotion-validity note from column R, and move that note
# text onto S22 (replacing the "licensed trade" note that used to live
# there).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Department (column C) per row, grouped by course type.
$businessRows = 2..12
$itRows = 13..15
$buildingRows = 16..16
$packageRows = 17..22

foreach ($r in $businessRows) {
    $ws.Range("C$r").Value = "Business"
}
foreach ($r in $itRows) {
    $ws.Range("C$r").Value = "Information Technology"
}
foreach ($r in $buildingRows) {
    $ws.Range("C$r").Value = "Building and Construction"
}
foreach ($r in $packageRows) {
    $ws.Range("C$r").Value = "Packages"
}

# Clear the old promotion-validity text from column R for every data row.
for ($r = 2; $r -le 22; $r++) {
    $ws.Range("R$r").Value = ""
}

# S22 now carries the promotion-validity note instead of the licensed
# trade disclaimer.
$ws.Range("S22").Value = "Promotion valid until  31th Dec 2021"
